$d = $word.ActiveDocument

# 1. The "_GoBack" bookmark currently sits on the "Creazione nuova categoria"
#    Heading2 paragraph. It needs to move to the edit that was actually made
#    below, so drop it from its current location first (a document can only
#    have one bookmark with a given name).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Fill in the Postcondizioni cell of the first table ("Elenco categorie
#    presenti nel sistema"), which currently just holds a placeholder "-".
$table = $d.Tables.Item(1)
$cell = $table.Cell(4, 2)
$cell.Range.Text = "Le Categorie così recuperate non hanno l’elenco degli Allergenti popolato"

# 3. Word leaves "_GoBack" marking the most recent edit, so re-create it
#    collapsed right after the text we just typed.
$editedCell = $table.Cell(4, 2)
$editedRange = $editedCell.Range
$editedRange.MoveEnd(1, -1) | Out-Null
$editedRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $editedRange) | Out-Null
